$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet tab and title to reflect new "through" date
$ws.Name = "Through 2022-12-08"

# Update header label in I1 (shared string "2022 (through 12-07)" -> "2022 (through 12-08)")
$ws.Range("I1").Value = "2022 (through 12-08)"

# Update November (row 12) and December (row 13) counts for the 2022 column (I),
# and recompute the Total (row 14)
$ws.Range("I12").Value = 118
$ws.Range("I13").Value = 34
$ws.Range("I14").Value = 1550
